$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90

$ws.Cells.Item($row, 1).Value = "Kindergarden"
$ws.Cells.Item($row, 2).Value = "Kindergarden Delft Reinier de Graaf"
$ws.Cells.Item($row, 3).Value = "KDV"

# The report date column stores plain text dates (e.g. "2024-09-23") in this
# workbook, not real Excel date serials. Force the cell to text first so the
# COM layer doesn't auto-convert the string into a date value, then strip the
# number-format override again so the cell keeps the same "no style" shape as
# every other data row.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-08-26"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
